$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''66.478.94'
$ws.Range('E2').Value = '''  +0.33%  '
$ws.Range('D3').Value = '''3.217.68'
$ws.Range('E3').Value = '''  +1.20%  '
$ws.Range('D4').Value = '''0.999'
$ws.Range('E4').Value = '''  -0.10%  '
$ws.Range('D5').Value = '''608.99'
$ws.Range('E5').Value = '''  +2.34%  '
$ws.Range('D6').Value = '''158.13'
$ws.Range('E6').Value = '''  +2.90%  '
$ws.Range('E7').Value = '''  -0.01%  '
$ws.Range('D8').Value = '''3.216.84'
$ws.Range('E8').Value = '''  +1.31%  '
$ws.Range('D9').Value = '''0.553'
$ws.Range('E9').Value = '''  +0.49%  '
$ws.Range('E10').Value = '''  +1.05%  '
$ws.Range('E11').Value = '''  -4.13%  '
$ws.Range('D12').Value = '''0.504'
$ws.Range('E12').Value = '''  -2.16%  '
$ws.Range('D13').Value = '''0.0000271'
$ws.Range('E13').Value = '''  +1.27%  '
$ws.Range('D14').Value = '''38.73'
$ws.Range('E14').Value = '''  -0.65%  '
$ws.Range('D15').Value = '''3.745.85'
$ws.Range('E15').Value = '''  +1.18%  '
$ws.Range('D16').Value = '''66.563.72'
$ws.Range('E16').Value = '''  +0.48%  '
$ws.Range('D17').Value = '''7.37'
$ws.Range('E17').Value = '''  -0.92%  '
$ws.Range('D18').Value = '''3.215.91'
$ws.Range('E18').Value = '''  +1.02%  '
$ws.Range('E19').Value = '''  +1.28%  '
$ws.Range('D20').Value = '''507.18'
$ws.Range('E20').Value = '''  -1.34%  '
$ws.Range('D21').Value = '''15.18'
$ws.Range('E21').Value = '''  -1.01%  '
$ws.Range('D22').Value = '''0.734'
$ws.Range('E22').Value = '''  -0.29%  '
$ws.Range('D23').Value = '''8.01'
$ws.Range('E23').Value = '''  -0.24%  '
$ws.Range('D24').Value = '''14.64'
$ws.Range('E24').Value = '''  -1.70%  '
$ws.Range('D25').Value = '''84.90'
$ws.Range('E25').Value = '''  -0.66%  '
$ws.Range('E26').Value = '''  +0.24%  '
$ws.Range('E27').Value = '''  +0.29%  '
$ws.Range('D28').Value = '''9.13'
$ws.Range('E28').Value = '''  -0.93%  '
$ws.Range('E29').Value = '''  +1.37%  '
$ws.Range('D30').Value = '''0.123'
$ws.Range('E30').Value = '''  +37.54%  '
$ws.Range('B31').Value = '''NEARProtocol'
$ws.Range('C31').Value = '''https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D31').Value = '''7.04'
$ws.Range('E31').Value = '''  -1.03%  '
$ws.Range('B32').Value = '''Stacks'
$ws.Range('C32').Value = '''https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D32').Value = '''2.94'
$ws.Range('E32').Value = '''  +0.65%  '
$ws.Range('D33').Value = '''28.18'
$ws.Range('E33').Value = '''  -0.16%  '
$ws.Range('E34').Value = '''  +0.12%  '
$ws.Range('E35').Value = '''  -4.56%  '
$ws.Range('D36').Value = '''6.49'
$ws.Range('E36').Value = '''  +0.00%  '
$ws.Range('D37').Value = '''503.94'
$ws.Range('E38').Value = '''  +1.09%  '
$ws.Range('D39').Value = '''0.0₃0774'
$ws.Range('E39').Value = '''  +14.89%  '
$ws.Range('B40').Value = '''dogwifhat'
$ws.Range('C40').Value = '''https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D40').Value = '''3.07'
$ws.Range('E40').Value = '''  +7.14%  '
$ws.Range('B41').Value = '''Kaspa'
$ws.Range('C41').Value = '''https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').Value = '''0.131'
$ws.Range('E41').Value = '''  +4.96%  '
$ws.Range('B42').Value = '''VeChain'
$ws.Range('C42').Value = '''https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').Value = '''0.0421'
$ws.Range('E42').Value = '''  -0.74%  '
$ws.Range('D43').Value = '''8.74'
$ws.Range('E43').Value = '''  -1.51%  '
$ws.Range('D44').Value = '''0.297'
$ws.Range('E44').Value = '''  -1.63%  '
$ws.Range('D45').Value = '''2.46'
$ws.Range('E45').Value = '''  +0.63%  '
$ws.Range('D46').Value = '''2.907.47'
$ws.Range('E46').Value = '''  +0.21%  '
$ws.Range('D47').Value = '''28.23'
$ws.Range('E47').Value = '''  -1.01%  '
$ws.Range('E48').Value = '''  +3.68%  '
$ws.Range('E49').Value = '''  -0.06%  '
$ws.Range('E50').Value = '''  -0.61%  '
$ws.Range('D51').Value = '''122.66'
$ws.Range('E51').Value = '''  -0.85%  '
